# QA updates to the US Core IG
# Replace comma-separated FHIR element/resource lists in column H with
# pipe-separated lists (and tighten one "Coverage.status + Coverage.period"
# entry to "Coverage.status+Coverage.period").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("uscdi-table.csv")

$updates = @{
    "H8"   = "Practitioner.name | RelatedPerson.name"
    "H10"  = "PractitionerRole.address | Practitioner.address | RelatedPerson.address"
    "H11"  = "PractitionerRole.telecom | Practitioner.telecom | RelatedPerson.telecom"
    "H13"  = "DocumentReference | DiagnosticReport"
    "H17"  = "DocumentReference | DiagnosticReport"
    "H18"  = "DocumentReference | DiagnosticReport"
    "H19"  = "DocumentReference | DiagnosticReport"
    "H20"  = "DocumentReference | DiagnosticReport"
    "H22"  = "Observation | DiagnosticReport"
    "H23"  = "Observation | DiagnosticReport"
    "H24"  = "Observation | DiagnosticReport"
    "H25"  = "Observation | DiagnosticReport"
    "H26"  = "Observation | DiagnosticReport"
    "H27"  = "Observation | DiagnosticReport"
    "H45"  = "Coverage.status+Coverage.period"
    "H47"  = "Coverage | relationship"
    "H52"  = "Observation | Condition | QuestionnaireResponse"
    "H58"  = "Observation | Condition | QuestionnaireResponse"
    "H59"  = "Observation | Condition | QuestionnaireResponse"
    "H60"  = "Observation | Condition | QuestionnaireResponse"
    "H61"  = "Observation | Condition | QuestionnaireResponse"
    "H65"  = "Observation | DiagnosticReport | Specimen"
    "H66"  = "Observation.code | DiagnosticReport.code"
    "H67"  = "Observation.value[x] | DiagnosticReport.result"
    "H68"  = "Observation.specimen | Specimen"
    "H69"  = "Observation.status | DiagnosticReport.status"
    "H78"  = "MedicationRequest | MedicationDispense | Medication"
    "H79"  = "Medication.code | MedicationRequest.medication[x]"
    "H86"  = "Patient | Observation"
    "H113" = "Condition.extension.where(url='http://hl7.org/fhir/StructureDefinition/condition-assertedDate') | Condition.onset[x] | Condition.onsetDate"
    "H115" = "Procedure | ServiceRequest"
    "H117" = "Procedure.pertormed[x] | DiagnosticReport.effective[x] | Immunization.occurrence[x]"
    "H118" = "ServiceRequest.reasonCode | ServiceRequest.reasonReference | Procedure.reasonCode | Procedure.reasonReference"
    "H119" = "ServiceRequest | Procedure"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
